$d = $word.ActiveDocument

# Update the four table cell values per the diff
$d.Content.Find.Execute("93 (30.5)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "93 (30.4)", 2)

$d.Content.Find.Execute("94 (30.8)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "94 (30.7)", 2)

$d.Content.Find.Execute("90 (29.5)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "91 (29.7)", 2)

$d.Content.Find.Execute("305 (100.0)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "306 (100.0)", 2)
